$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H29").Value = 1999
$ws.Range("J29").Value = 2094.95
$ws.Range("L29").Value = 6284.849999999999
$ws.Range("N29").Value = -6846.849999999999
$ws.Range("H64").Value = 3718.8125
$ws.Range("I64").Value = 3540.1
$ws.Range("K64").Value = 3540.1
$ws.Range("M64").Value = -3292.1
$ws.Range("H67").Value = 3718.8125
$ws.Range("I67").Value = 3540.1
$ws.Range("K67").Value = 3540.1
$ws.Range("M67").Value = -2682.1
$ws.Range("H74").Value = 4299.375
$ws.Range("I74").Value = 3998.5
$ws.Range("J74").Value = 4399.6665
$ws.Range("K74").Value = 3998.5
$ws.Range("L74").Value = 4399.6665
$ws.Range("M74").Value = -3062.5
$ws.Range("N74").Value = -6271.6665
$ws.Range("H76").Value = 3042.8215
$ws.Range("I76").Value = 3031.96
$ws.Range("K76").Value = 3031.96
$ws.Range("M76").Value = -2716.96
$ws.Range("H77").Value = 4299.375
$ws.Range("I77").Value = 3998.5
$ws.Range("J77").Value = 4399.6665
$ws.Range("K77").Value = 19992.5
$ws.Range("L77").Value = 21998.3325
$ws.Range("M77").Value = -15312.5
$ws.Range("N77").Value = -31358.3325
$ws.Range("H79").Value = 3042.8215
$ws.Range("I79").Value = 3031.96
$ws.Range("K79").Value = 3031.96
$ws.Range("M79").Value = -1939.96
$ws.Range("H88").Value = 2939.0356
$ws.Range("I88").Value = 1952.4546
$ws.Range("K88").Value = 1952.4546
$ws.Range("M88").Value = -1546.4546
$ws.Range("H91").Value = 2939.0356
$ws.Range("I91").Value = 1952.4546
$ws.Range("K91").Value = 1952.4546
$ws.Range("M91").Value = -548.4546
$ws.Range("H96").Value = 1665.4
$ws.Range("I96").Value = 978.8570999999999
$ws.Range("J96").Value = 2266.125
$ws.Range("K96").Value = 2936.5713
$ws.Range("L96").Value = 6798.375
$ws.Range("M96").Value = -1563.5713
$ws.Range("N96").Value = -9544.375
$ws.Range("H106").Value = 6027.8667
$ws.Range("I106").Value = 5316.6665
$ws.Range("J106").Value = 6502
$ws.Range("K106").Value = 5316.6665
$ws.Range("L106").Value = 6502
$ws.Range("M106").Value = -4685.6665
$ws.Range("N106").Value = -7764
$ws.Range("H124").Value = 53333.332
$ws.Range("J124").Value = 53333.332
$ws.Range("L124").Value = 53333.332
$ws.Range("N124").Value = -63153.332
$ws.Range("H127").Value = 52632532
$ws.Range("I127").Value = 66667470
$ws.Range("J127").Value = 1516.5
$ws.Range("K127").Value = 200002410
$ws.Range("L127").Value = 4549.5
$ws.Range("M127").Value = -199997450
$ws.Range("N127").Value = -14469.5
$ws.Range("H137").Value = 2605.394
$ws.Range("I137").Value = 1916.4348
$ws.Range("J137").Value = 4190
$ws.Range("K137").Value = 5749.3044
$ws.Range("L137").Value = 12570
$ws.Range("M137").Value = -3199.3044
$ws.Range("N137").Value = -17670

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H25").Value = 8583.333000000001
$ws.Range("I25").Value = 2000
$ws.Range("J25").Value = 9900
$ws.Range("K25").Value = 2000
$ws.Range("L25").Value = 9900
$ws.Range("M25").Value = -1598
$ws.Range("N25").Value = -10704
$ws.Range("H32").Value = 3189.3064
$ws.Range("I32").Value = 2682.8035
$ws.Range("K32").Value = 2682.8035
$ws.Range("M32").Value = -2395.8035
$ws.Range("H45").Value = 1473.7273
$ws.Range("I45").Value = 1019.5
$ws.Range("J45").Value = 2018.8
$ws.Range("K45").Value = 1019.5
$ws.Range("L45").Value = 2018.8
$ws.Range("M45").Value = -642.5
$ws.Range("N45").Value = -2772.8
$ws.Range("H63").Value = 3732.5
$ws.Range("I63").Value = 3531.6667
$ws.Range("K63").Value = 3531.6667
$ws.Range("M63").Value = -2845.6667
$ws.Range("H66").Value = 3732.5
$ws.Range("I66").Value = 3531.6667
$ws.Range("K66").Value = 17658.3335
$ws.Range("M66").Value = -14226.3335
$ws.Range("H88").Value = 2367.476
$ws.Range("I88").Value = 2121.7
$ws.Range("K88").Value = 2121.7
$ws.Range("M88").Value = -1715.7
$ws.Range("H91").Value = 2367.476
$ws.Range("I91").Value = 2121.7
$ws.Range("K91").Value = 2121.7
$ws.Range("M91").Value = -717.6999999999998
$ws.Range("H92").Value = 28800
$ws.Range("J92").Value = 28800
$ws.Range("L92").Value = 28800
$ws.Range("N92").Value = -33792
$ws.Range("H122").Value = 2011.0769
$ws.Range("I122").Value = 1403.2858
$ws.Range("J122").Value = 2720.1667
$ws.Range("K122").Value = 4209.857400000001
$ws.Range("L122").Value = 8160.500100000001
$ws.Range("M122").Value = -1759.857400000001
$ws.Range("N122").Value = -13060.5001

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1764.6428
$ws.Range("I86").Value = 1642.0834
$ws.Range("K86").Value = 1642.0834
$ws.Range("M86").Value = -519.0834
$ws.Range("H89").Value = 1764.6428
$ws.Range("I89").Value = 1642.0834
$ws.Range("K89").Value = 8210.416999999999
$ws.Range("M89").Value = -2594.416999999999
$ws.Range("H92").Value = 71399.60000000001
$ws.Range("J92").Value = 71399.60000000001
$ws.Range("L92").Value = 71399.60000000001
$ws.Range("N92").Value = -76391.60000000001
$ws.Range("H105").Value = 2182.125
$ws.Range("I105").Value = 1806
$ws.Range("K105").Value = 1806
$ws.Range("M105").Value = -59
$ws.Range("H108").Value = 28974
$ws.Range("J108").Value = 28974
$ws.Range("L108").Value = 28974
$ws.Range("N108").Value = -36654

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 21740610
$ws.Range("I31").Value = 33334366
$ws.Range("J31").Value = 2319.3125
$ws.Range("K31").Value = 33334366
$ws.Range("L31").Value = 2319.3125
$ws.Range("M31").Value = -33334071
$ws.Range("N31").Value = -2909.3125
$ws.Range("H34").Value = 21740610
$ws.Range("I34").Value = 33334366
$ws.Range("J34").Value = 2319.3125
$ws.Range("K34").Value = 33334366
$ws.Range("L34").Value = 2319.3125
$ws.Range("M34").Value = -33334164
$ws.Range("N34").Value = -2723.3125
$ws.Range("H48").Value = 9966.333000000001
$ws.Range("J48").Value = 9966.333000000001
$ws.Range("L48").Value = 9966.333000000001
$ws.Range("N48").Value = -10918.333
$ws.Range("H62").Value = 3015.8462
$ws.Range("I62").Value = 2300
$ws.Range("K62").Value = 2300
$ws.Range("M62").Value = -1676
$ws.Range("H65").Value = 3015.8462
$ws.Range("I65").Value = 2300
$ws.Range("K65").Value = 11500
$ws.Range("M65").Value = -8380

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H93").Value = 6997
$ws.Range("J93").Value = 6997
$ws.Range("L93").Value = 20991
$ws.Range("N93").Value = -24735
$ws.Range("H113").Value = 526.1321
$ws.Range("I113").Value = 530.64514
$ws.Range("J113").Value = 519.7727
$ws.Range("K113").Value = 1591.93542
$ws.Range("L113").Value = 1559.3181
$ws.Range("M113").Value = 578.0645800000002
$ws.Range("N113").Value = -5899.3181
$ws.Range("H131").Value = 887.36365
$ws.Range("J131").Value = 891.337
$ws.Range("L131").Value = 2674.011
$ws.Range("N131").Value = -12754.011

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 4366.467
$ws.Range("I12").Value = 2625.75
$ws.Range("J12").Value = 4999.4546
$ws.Range("K12").Value = 2625.75
$ws.Range("L12").Value = 4999.4546
$ws.Range("M12").Value = -2485.75
$ws.Range("N12").Value = -5279.4546
$ws.Range("H24").Value = 12000
$ws.Range("I24").Value = 0
$ws.Range("J24").Value = 12000
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = 12000
$ws.Range("M24").ClearContents()
$ws.Range("N24").Value = -12346
$ws.Range("H70").Value = 4944.18
$ws.Range("I70").Value = 4869.4595
$ws.Range("J70").Value = 5059.375
$ws.Range("K70").Value = 4869.4595
$ws.Range("L70").Value = 5059.375
$ws.Range("M70").Value = -4599.4595
$ws.Range("N70").Value = -5599.375
$ws.Range("H73").Value = 4944.18
$ws.Range("I73").Value = 4869.4595
$ws.Range("J73").Value = 5059.375
$ws.Range("K73").Value = 4869.4595
$ws.Range("L73").Value = 5059.375
$ws.Range("M73").Value = -3933.4595
$ws.Range("N73").Value = -6931.375

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H20").Value = 37502.5
$ws.Range("I20").Value = 50005
$ws.Range("J20").Value = 25000
$ws.Range("K20").Value = 50005
$ws.Range("L20").Value = 25000
$ws.Range("M20").Value = -49779
$ws.Range("N20").Value = -25452

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 4005
$ws.Range("I22").Value = 1000
$ws.Range("K22").Value = 1000
$ws.Range("M22").Value = -707
$ws.Range("H80").Value = 30000
$ws.Range("J80").Value = 30000
$ws.Range("L80").Value = 30000
$ws.Range("N80").Value = -31996
$ws.Range("H83").Value = 30000
$ws.Range("J83").Value = 30000
$ws.Range("L83").Value = 90000
$ws.Range("N83").Value = -99984
$ws.Range("H128").Value = 34000
$ws.Range("J128").Value = 34000
$ws.Range("L128").Value = 34000
$ws.Range("N128").Value = -43960
$ws.Range("H135").Value = 39400
$ws.Range("J135").Value = 39400
$ws.Range("L135").Value = 39400
$ws.Range("N135").Value = -49540
